$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove K2 value (was "XX")
$ws.Range("K2").ClearContents()

# Add new header cells L1 and N1
$ws.Range("L1").Value = "MI"
$ws.Range("N1").Value = "VT spreadsheet"

# Add new cells in row 20: H20, I20, J20 replicate G20's value; K20 is a single space
$g20 = $ws.Range("G20").Value2
$ws.Range("H20").Value = $g20
$ws.Range("I20").Value = $g20
$ws.Range("J20").Value = $g20
$ws.Range("K20").Value = " "

# Update selection to match diff (active cell J10)
$ws.Range("J10").Select()
